# Updated the remaining MCH collections
# Adds two new catalogue rows (MCH338-1 and MCH338-2) to Sheet1, matching
# the columns established by the header row:
#   A identifier | B alternativeIdentifiers | C title | D date_s
#   E levelOfDescription | F extentAndMedium | G notes | H file_path

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 : MCH338-1 -------------------------------------------------
$ws.Range("A2").Value = "MCH338-1"
$ws.Range("C2").Value = "BIOGRAPHIES- ALEXANDER, NEVILLE. AUERBACH, FRANTZ. CAROLUS, CHERYL. DUNCAN, SHEENA. GHANDI, ELA. MANDELA, NELSON. MBEKI, GOVAN. MOGOBA, STANLEY. MOMPATI, RUTH. SLOVO, JOE."
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 33H | GRAP COUNT NUMER: NONE"

# --- Row 3 : MCH338-2 -------------------------------------------------
$ws.Range("A3").Value = "MCH338-2"
$ws.Range("C3").Value = "BIOGRAPHIES- NADINE GORDIMER, ITUMELENG MOSALA, TREVOR HUDDLESTONE, BEYERS NAUDE, ALBERTINA SISULU, EBRAHIM RASOOL, FRANK CHIKANE, FATIMA MEER, RAY ALEXANDER, NELSON MANDELA, CHRIS HANI."
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 33H | GRAP COUNT NUMER: NONE"

# D2/D3 (date_s) and H2/H3 (file_path) are left blank on purpose, but they
# still pick up the new row formatting below (empty-but-styled cells), so
# touch them too.

# --- Formatting: new rows use 10pt Calibri, automatic (theme) text colour,
# matching the rest of the data area. Work one contiguous block at a time
# (per column-group) so every cell in that block is re-styled together.
$columnGroups = @("A2:A3", "C2:D3", "E2:E3", "G2:H3")
foreach ($addr in $columnGroups) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Calibri"
    $rng.Font.ThemeColor = 1
}

# extentAndMedium (F) carries its own (slightly distinct) style in the
# source workbook, so it is styled as its own block.
$fRng = $ws.Range("F2:F3")
$fRng.Font.Name = "Calibri"
$fRng.Font.ThemeColor = 1
